# Refresh cached market-price / profit figures for several Leve rows across
# multiple crafting-job sheets (scheduled market-data pull). Pure value
# updates -- no formulas, formatting, or structural changes involved.
$wb = $excel.ActiveWorkbook

# ---- ALC sheet ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3818873.8
$ws.Range("I40").Value = 7815050
$ws.Range("J40").Value = 912563.8
$ws.Range("K40").Value = 7815050
$ws.Range("L40").Value = 912563.8
$ws.Range("M40").Value = -7814875
$ws.Range("N40").Value = -912913.8

$ws.Range("H51").Value = 20309.273
$ws.Range("I51").Value = 41130.4
$ws.Range("K51").Value = 41130.4
$ws.Range("M51").Value = -40646.4

$ws.Range("H64").Value = 3680
$ws.Range("I64").Value = 3542.8572
$ws.Range("J64").Value = 4000
$ws.Range("K64").Value = 3542.8572
$ws.Range("L64").Value = 4000
$ws.Range("M64").Value = -3294.8572
$ws.Range("N64").Value = -4496

$ws.Range("H67").Value = 3680
$ws.Range("I67").Value = 3542.8572
$ws.Range("J67").Value = 4000
$ws.Range("K67").Value = 3542.8572
$ws.Range("L67").Value = 4000
$ws.Range("M67").Value = -2684.8572
$ws.Range("N67").Value = -5716

$ws.Range("H107").Value = 84207.914
$ws.Range("I107").Value = 143330.72
$ws.Range("J107").Value = 1436
$ws.Range("K107").Value = 143330.72
$ws.Range("L107").Value = 1436
$ws.Range("M107").Value = -141410.72
$ws.Range("N107").Value = -5276

$ws.Range("H121").Value = 883.44446
$ws.Range("J121").Value = 1021.5714
$ws.Range("L121").Value = 3064.7142
$ws.Range("N121").Value = -6558.7142

$ws.Range("H131").Value = 1527.8889
$ws.Range("I131").Value = 458.39285
$ws.Range("K131").Value = 1375.17855
$ws.Range("M131").Value = 3664.82145

$ws.Range("H137").Value = 1160.4032
$ws.Range("I137").Value = 958.1852
$ws.Range("K137").Value = 2874.5556
$ws.Range("M137").Value = -324.5556000000001

$ws.Range("H141").Value = 1937.2858
$ws.Range("I141").Value = 1093.5
$ws.Range("J141").Value = 7000
$ws.Range("K141").Value = 3280.5
$ws.Range("L141").Value = 21000
$ws.Range("M141").Value = 1899.5
$ws.Range("N141").Value = -31360

# ---- ARM sheet ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11239661
$ws.Range("I32").Value = 3398.4824
$ws.Range("K32").Value = 3398.4824
$ws.Range("M32").Value = -3111.4824

$ws.Range("H63").Value = 2304.9092
$ws.Range("I63").Value = 2310.4
$ws.Range("J63").Value = 2250
$ws.Range("K63").Value = 2310.4
$ws.Range("L63").Value = 2250
$ws.Range("M63").Value = -1624.4
$ws.Range("N63").Value = -3622

$ws.Range("H66").Value = 2304.9092
$ws.Range("I66").Value = 2310.4
$ws.Range("J66").Value = 2250
$ws.Range("K66").Value = 11552
$ws.Range("L66").Value = 11250
$ws.Range("M66").Value = -8120
$ws.Range("N66").Value = -18114

# ---- BSM sheet ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 4504637.5
$ws.Range("I22").Value = 4504637.5
$ws.Range("K22").Value = 4504637.5
$ws.Range("M22").Value = -4504464.5

$ws.Range("H82").Value = 5995
$ws.Range("I82").Value = 5995
$ws.Range("K82").Value = 5995
$ws.Range("M82").Value = -5612

$ws.Range("H85").Value = 5995
$ws.Range("I85").Value = 5995
$ws.Range("K85").Value = 5995
$ws.Range("M85").Value = -4669

# ---- CRP sheet ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1258.8572
$ws.Range("I31").Value = 843.36365
$ws.Range("K31").Value = 843.36365
$ws.Range("M31").Value = -548.36365

$ws.Range("M32").Value = -4080.6665
$ws.Range("N32").Value = -10632  # new cell added in target

$ws.Range("H34").Value = 1258.8572
$ws.Range("I34").Value = 843.36365
$ws.Range("K34").Value = 843.36365
$ws.Range("M34").Value = -641.36365

$ws.Range("H58").Value = 47619976
$ws.Range("I58").Value = 125001090
$ws.Range("J58").Value = 832.9231
$ws.Range("K58").Value = 125001090
$ws.Range("L58").Value = 832.9231
$ws.Range("M58").Value = -125000887
$ws.Range("N58").Value = -1238.9231

$ws.Range("H136").Value = 47619976
$ws.Range("I136").Value = 125001090
$ws.Range("J136").Value = 832.9231
$ws.Range("K136").Value = 375003270
$ws.Range("L136").Value = 2498.7693
$ws.Range("M136").Value = -375000720
$ws.Range("N136").Value = -7598.7693

# ---- CUL sheet ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 24695500
$ws.Range("J5").Value = 17334.166
$ws.Range("L5").Value = 52002.49800000001
$ws.Range("N5").Value = -52226.49800000001

$ws.Range("L81").Value = 9611.25
$ws.Range("M81").ClearContents()  # cell removed in target

$ws.Range("L84").Value = 28833.75
$ws.Range("M84").ClearContents()  # cell removed in target

$ws.Range("H135").Value = 24695500
$ws.Range("J135").Value = 17334.166
$ws.Range("L135").Value = 156007.494
$ws.Range("N135").Value = -161077.494

# ---- GSM sheet ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1380.3572
$ws.Range("I113").Value = 1150
$ws.Range("J113").Value = 1472.5
$ws.Range("K113").Value = 1150
$ws.Range("L113").Value = 1472.5
$ws.Range("M113").Value = 1020
$ws.Range("N113").Value = -5812.5

$ws.Range("H132").Value = 4123.8696
$ws.Range("I132").Value = 1546.5333
$ws.Range("J132").Value = 8956.375
$ws.Range("K132").Value = 4639.5999
$ws.Range("L132").Value = 26869.125
$ws.Range("M132").Value = -2109.5999
$ws.Range("N132").Value = -31929.125

# ---- LTW sheet ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 603876.2
$ws.Range("I22").Value = 974089.25
$ws.Range("J22").Value = 2280
$ws.Range("K22").Value = 974089.25
$ws.Range("L22").Value = 2280
$ws.Range("M22").Value = -973794.25
$ws.Range("N22").Value = -2870

$ws.Range("H27").Value = 603876.2
$ws.Range("I27").Value = 974089.25
$ws.Range("J27").Value = 2280
$ws.Range("K27").Value = 974089.25
$ws.Range("L27").Value = 2280
$ws.Range("M27").Value = -973982.25
$ws.Range("N27").Value = -2494
